$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 169, pushing existing rows 169-174 down to 170-175.
$ws.Rows.Item(169).Insert()

# Populate the new row 169 with the new weekly price record (copy of row 170's
# style for column D, same as the rest of this data block).
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169, 3).Value = "Bíobío"
$ws.Cells.Item(169, 4).Value = 44714
$ws.Cells.Item(169, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 100112003
$ws.Cells.Item(169, 7).Value = "Ajo"
$ws.Cells.Item(169, 8).Value = "Chino"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 220
$ws.Cells.Item(169, 11).Value = 16000
$ws.Cells.Item(169, 12).Value = 17000
$ws.Cells.Item(169, 13).Value = 16455
$ws.Cells.Item(169, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(169, 15).Value = "China"
$ws.Cells.Item(169, 16).Value = 1646
$ws.Cells.Item(169, 17).Value = 10
$ws.Cells.Item(169, 18).Value = "Hortaliza"
